# Updates cryptos list prices/volumes (and re-ranks two coin pairs) per the
# Mon Apr 22 16:15:22 UTC 2024 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'66.077.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = "'3.188.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.38%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = "'594.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.66%  '
$ws.Range("D6").Value = "'154.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.50%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.27%  '
$ws.Range("D8").Value = "'3.189.95"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.60%  '
$ws.Range("D9").Value = "'0.541"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.02%  '
$ws.Range("D10").Value = "'0.159"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.35%  '
$ws.Range("D11").Value = "'6.06"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.11%  '
$ws.Range("D12").Value = "'0.517"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.31%  '
$ws.Range("E13").Value = '  +1.04%  '
$ws.Range("D14").Value = "'39.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.16%  '
$ws.Range("D15").Value = "'3.709.99"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = "'7.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.52%  '
$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").Value = "'65.977.79"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = "'3.186.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.26%  '
$ws.Range("D19").Value = "'0.112"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("D20").Value = "'510.44"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.66%  '
$ws.Range("D21").Value = "'15.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.96%  '
$ws.Range("D22").Value = "'0.740"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.71%  '
$ws.Range("D23").Value = "'8.05"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.57%  '
$ws.Range("D24").Value = "'15.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.75%  '
$ws.Range("D25").Value = "'85.01"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.18%  '
$ws.Range("D26").Value = "'0.997"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").Value = "'9.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.75%  '
$ws.Range("D28").Value = "'2.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.24%  '
$ws.Range("D29").Value = "'2.31"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.97%  '
$ws.Range("D30").Value = "'7.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +14.63%  '
$ws.Range("D31").Value = "'2.89"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.42%  '
$ws.Range("D32").Value = "'28.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("E33").Value = '  +3.07%  '
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = "'6.50"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.44%  '
$ws.Range("D36").Value = "'54.80"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.06%  '
$ws.Range("E37").Value = '  +0.42%  '
$ws.Range("D38").Value = "'482.96"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.42%  '
$ws.Range("D39").Value = "'0.0421"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.05%  '
$ws.Range("D40").Value = "'8.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.19%  '
$ws.Range("D41").Value = "'0.303"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +7.43%  '
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = "'2.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.17%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").Value = "'0.122"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.18%  '
$ws.Range("D44").Value = "'0.0₃0658"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +12.71%  '
$ws.Range("D45").Value = "'2.909.87"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.16%  '
$ws.Range("D46").Value = "'2.42"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("D47").Value = "'28.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.76%  '
$ws.Range("D48").Value = "'0.999"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.02%  '
$ws.Range("E49").Value = '  +2.38%  '
$ws.Range("D50").Value = "'2.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.78%  '
$ws.Range("D51").Value = "'2.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.13%  '
